$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2.0
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.352356
$ws.Range("H2").Value = 4.057068
$ws.Range("I2").Value = 0.08465764667810434
$ws.Range("J2").Value = 0.1020450778651008
$ws.Range("M2").Value = 55.783591
$ws.Range("N2").Value = 167.350773
$ws.Range("O2").Value = 0.2332214199005771
$ws.Range("P2").Value = 0.2394371967339281
$ws.Range("Q2").Value = 75.439273990396
$ws.Range("R2").Value = 678.9534659135641
$ws.Range("S2").Value = 0.01974397656370887
$ws.Range("T2").Value = 0.02443338738451515
$ws.Range("E3").Value = 2.0
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.352356
$ws.Range("H3").Value = 4.057068
$ws.Range("I3").Value = 0.08465764667810434
$ws.Range("J3").Value = 0.1020450778651008
$ws.Range("O3").Value = 0.5297503589663128
$ws.Range("P3").Value = 0.5438691736537713
$ws.Range("Q3").Value = 171.3563981113187
$ws.Range("R3").Value = 1542.207583001868
$ws.Range("S3").Value = 0.04484741871696905
$ws.Range("T3").Value = 0.05549917217392714
$ws.Range("E4").Value = 2.0
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.352356
$ws.Range("H4").Value = 4.057068
$ws.Range("I4").Value = 0.08465764667810434
$ws.Range("J4").Value = 0.1020450778651008
$ws.Range("M4").Value = 23.03749833333333
$ws.Range("N4").Value = 69.112495
$ws.Range("O4").Value = 0.09631574403765399
$ws.Range("P4").Value = 0.09888273454277752
$ws.Range("Q4").Value = 31.15489909607333
$ws.Range("R4").Value = 280.39409186466
$ws.Range("S4").Value = 0.008153864228278445
$ws.Range("T4").Value = 0.01009049634593183
$ws.Range("E5").Value = 2.0
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 1.352356
$ws.Range("H5").Value = 4.057068
$ws.Range("I5").Value = 0.08465764667810434
$ws.Range("J5").Value = 0.1020450778651008
$ws.Range("M5").Value = 18.627865
$ws.Range("N5").Value = 37.25573
$ws.Range("O5").Value = 0.07787983970082285
$ws.Range("P5").Value = 0.05330365312071852
$ws.Range("Q5").Value = 25.19150499994
$ws.Range("R5").Value = 151.14902999964
$ws.Range("S5").Value = 0.006593123952739664
$ws.Range("T5").Value = 0.005439375433198046
$ws.Range("E6").Value = 2.0
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 1.352356
$ws.Range("H6").Value = 4.057068
$ws.Range("I6").Value = 0.08465764667810434
$ws.Range("J6").Value = 0.1020450778651008
$ws.Range("M6").Value = 15.028766
$ws.Range("N6").Value = 45.086298
$ws.Range("O6").Value = 0.06283263739463307
$ws.Range("P6").Value = 0.06450724194880479
$ws.Range("Q6").Value = 20.324241872696
$ws.Range("R6").Value = 182.918176854264
$ws.Range("S6").Value = 0.005319263216408293
$ws.Range("T6").Value = 0.006582646527528683
$ws.Range("I7").Value = 0.2001834746890973
$ws.Range("J7").Value = 0.2412982059332244
$ws.Range("M7").Value = 55.783591
$ws.Range("N7").Value = 167.350773
$ws.Range("O7").Value = 0.2332214199005771
$ws.Range("P7").Value = 0.2394371967339281
$ws.Range("Q7").Value = 178.385492486483
$ws.Range("R7").Value = 1605.469432378347
$ws.Range("S7").Value = 0.04668707420762251
$ws.Range("T7").Value = 0.05777576600557735
$ws.Range("I8").Value = 0.2001834746890973
$ws.Range("J8").Value = 0.2412982059332244
$ws.Range("O8").Value = 0.5297503589663128
$ws.Range("P8").Value = 0.5438691736537713
$ws.Range("S8").Value = 0.1060472675756731
$ws.Range("T8").Value = 0.1312346558650403
$ws.Range("I9").Value = 0.2001834746890973
$ws.Range("J9").Value = 0.2412982059332244
$ws.Range("M9").Value = 23.03749833333333
$ws.Range("N9").Value = 69.112495
$ws.Range("O9").Value = 0.09631574403765399
$ws.Range("P9").Value = 0.09888273454277752
$ws.Range("Q9").Value = 73.66961165781166
$ws.Range("R9").Value = 663.026504920305
$ws.Range("S9").Value = 0.01928082030872328
$ws.Range("T9").Value = 0.02386022644294349
$ws.Range("I10").Value = 0.2001834746890973
$ws.Range("J10").Value = 0.2412982059332244
$ws.Range("M10").Value = 18.627865
$ws.Range("N10").Value = 37.25573
$ws.Range("O10").Value = 0.07787983970082285
$ws.Range("P10").Value = 0.05330365312071852
$ws.Range("Q10").Value = 59.568428859245
$ws.Range("R10").Value = 357.41057315547
$ws.Range("S10").Value = 0.01559025691954062
$ws.Range("T10").Value = 0.0128620758677163
$ws.Range("I11").Value = 0.2001834746890973
$ws.Range("J11").Value = 0.2412982059332244
$ws.Range("M11").Value = 15.028766
$ws.Range("N11").Value = 45.086298
$ws.Range("O11").Value = 0.06283263739463307
$ws.Range("P11").Value = 0.06450724194880479
$ws.Range("Q11").Value = 48.059183288758
$ws.Range("R11").Value = 432.532649598822
$ws.Range("S11").Value = 0.01257805567753776
$ws.Range("T11").Value = 0.01556548175194703
$ws.Range("G12").Value = 2.438989
$ws.Range("H12").Value = 7.316967
$ws.Range("I12").Value = 0.1526810019061423
$ws.Range("J12").Value = 0.1840394263175705
$ws.Range("M12").Value = 55.783591
$ws.Range("N12").Value = 167.350773
$ws.Range("O12").Value = 0.2332214199005771
$ws.Range("P12").Value = 0.2394371967339281
$ws.Range("Q12").Value = 136.055564829499
$ws.Range("R12").Value = 1224.500083465491
$ws.Range("S12").Value = 0.03560848005639324
$ws.Range("T12").Value = 0.04406588432599939
$ws.Range("G13").Value = 2.438989
$ws.Range("H13").Value = 7.316967
$ws.Range("I13").Value = 0.1526810019061423
$ws.Range("J13").Value = 0.1840394263175705
$ws.Range("O13").Value = 0.5297503589663128
$ws.Range("P13").Value = 0.5438691736537713
$ws.Range("Q13").Value = 309.0431588081296
$ws.Range("R13").Value = 2781.388429273167
$ws.Range("S13").Value = 0.08088281556711518
$ws.Range("T13").Value = 0.1000933707110512
$ws.Range("G14").Value = 2.438989
$ws.Range("H14").Value = 7.316967
$ws.Range("I14").Value = 0.1526810019061423
$ws.Range("J14").Value = 0.1840394263175705
$ws.Range("M14").Value = 23.03749833333333
$ws.Range("N14").Value = 69.112495
$ws.Range("O14").Value = 0.09631574403765399
$ws.Range("P14").Value = 0.09888273454277752
$ws.Range("Q14").Value = 56.18820502251832
$ws.Range("R14").Value = 505.693845202665
$ws.Range("S14").Value = 0.01470558429900456
$ws.Range("T14").Value = 0.01819832173796539
$ws.Range("G15").Value = 2.438989
$ws.Range("H15").Value = 7.316967
$ws.Range("I15").Value = 0.1526810019061423
$ws.Range("J15").Value = 0.1840394263175705
$ws.Range("M15").Value = 18.627865
$ws.Range("N15").Value = 37.25573
$ws.Range("O15").Value = 0.07787983970082285
$ws.Range("P15").Value = 0.05330365312071852
$ws.Range("Q15").Value = 45.43315782848499
$ws.Range("R15").Value = 272.59894697091
$ws.Range("S15").Value = 0.01189077195381139
$ws.Range("T15").Value = 0.009809973740967813
$ws.Range("G16").Value = 2.438989
$ws.Range("H16").Value = 7.316967
$ws.Range("I16").Value = 0.1526810019061423
$ws.Range("J16").Value = 0.1840394263175705
$ws.Range("M16").Value = 15.028766
$ws.Range("N16").Value = 45.086298
$ws.Range("O16").Value = 0.06283263739463307
$ws.Range("P16").Value = 0.06450724194880479
$ws.Range("Q16").Value = 36.654994957574
$ws.Range("R16").Value = 329.894954618166
$ws.Range("S16").Value = 0.009593350029817921
$ws.Range("T16").Value = 0.01187187580158675
$ws.Range("G17").Value = 8.1656255
$ws.Range("H17").Value = 16.331251
$ws.Range("I17").Value = 0.5111691289015016
$ws.Range("J17").Value = 0.4107704825084287
$ws.Range("M17").Value = 55.783591
$ws.Range("N17").Value = 167.350773
$ws.Range("O17").Value = 0.2332214199005771
$ws.Range("P17").Value = 0.2394371967339281
$ws.Range("Q17").Value = 455.5079131511706
$ws.Range("R17").Value = 2733.047478907024
$ws.Range("S17").Value = 0.1192155900517493
$ws.Range("T17").Value = 0.09835373283286121
$ws.Range("G18").Value = 8.1656255
$ws.Range("H18").Value = 16.331251
$ws.Range("I18").Value = 0.5111691289015016
$ws.Range("J18").Value = 0.4107704825084287
$ws.Range("O18").Value = 0.5297503589663128
$ws.Range("P18").Value = 0.5438691736537713
$ws.Range("Q18").Value = 1034.662599201642
$ws.Range("R18").Value = 6207.975595209852
$ws.Range("S18").Value = 0.2707920295280679
$ws.Range("T18").Value = 0.2234054028832201
$ws.Range("G19").Value = 8.1656255
$ws.Range("H19").Value = 16.331251
$ws.Range("I19").Value = 0.5111691289015016
$ws.Range("J19").Value = 0.4107704825084287
$ws.Range("M19").Value = 23.03749833333333
$ws.Range("N19").Value = 69.112495
$ws.Range("O19").Value = 0.09631574403765399
$ws.Range("P19").Value = 0.09888273454277752
$ws.Range("Q19").Value = 188.1155838468742
$ws.Range("R19").Value = 1128.693503081245
$ws.Range("S19").Value = 0.04923363497922758
$ws.Range("T19").Value = 0.04061810857988959
$ws.Range("G20").Value = 8.1656255
$ws.Range("H20").Value = 16.331251
$ws.Range("I20").Value = 0.5111691289015016
$ws.Range("J20").Value = 0.4107704825084287
$ws.Range("M20").Value = 18.627865
$ws.Range("N20").Value = 37.25573
$ws.Range("O20").Value = 0.07787983970082285
$ws.Range("P20").Value = 0.05330365312071852
$ws.Range("Q20").Value = 152.1081694545575
$ws.Range("R20").Value = 608.4326778182301
$ws.Range("S20").Value = 0.0398097698188582
$ws.Range("T20").Value = 0.02189556731185946
$ws.Range("G21").Value = 8.1656255
$ws.Range("H21").Value = 16.331251
$ws.Range("I21").Value = 0.5111691289015016
$ws.Range("J21").Value = 0.4107704825084287
$ws.Range("M21").Value = 15.028766
$ws.Range("N21").Value = 45.086298
$ws.Range("O21").Value = 0.06283263739463307
$ws.Range("P21").Value = 0.06450724194880479
$ws.Range("Q21").Value = 122.719274883133
$ws.Range("R21").Value = 736.3156492987981
$ws.Range("S21").Value = 0.0321181045235985
$ws.Range("T21").Value = 0.0264976709005985
$ws.Range("G22").Value = 0.819627
$ws.Range("H22").Value = 2.458881
$ws.Range("I22").Value = 0.05130874782515448
$ws.Range("J22").Value = 0.06184680737567548
$ws.Range("M22").Value = 55.783591
$ws.Range("N22").Value = 167.350773
$ws.Range("O22").Value = 0.2332214199005771
$ws.Range("P22").Value = 0.2394371967339281
$ws.Range("Q22").Value = 45.721737340557
$ws.Range("R22").Value = 411.495636065013
$ws.Range("S22").Value = 0.01196629902110318
$ws.Range("T22").Value = 0.01480842618497496
$ws.Range("G23").Value = 0.819627
$ws.Range("H23").Value = 2.458881
$ws.Range("I23").Value = 0.05130874782515448
$ws.Range("J23").Value = 0.06184680737567548
$ws.Range("O23").Value = 0.5297503589663128
$ws.Range("P23").Value = 0.5438691736537713
$ws.Range("Q23").Value = 103.854554950609
$ws.Range("R23").Value = 934.690994555481
$ws.Range("S23").Value = 0.02718082757848761
$ws.Range("T23").Value = 0.03363657202053259
$ws.Range("G24").Value = 0.819627
$ws.Range("H24").Value = 2.458881
$ws.Range("I24").Value = 0.05130874782515448
$ws.Range("J24").Value = 0.06184680737567548
$ws.Range("M24").Value = 23.03749833333333
$ws.Range("N24").Value = 69.112495
$ws.Range("O24").Value = 0.09631574403765399
$ws.Range("P24").Value = 0.09888273454277752
$ws.Range("Q24").Value = 18.882155646455
$ws.Range("R24").Value = 169.939400818095
$ws.Range("S24").Value = 0.004941840222420115
$ws.Range("T24").Value = 0.006115581436047213
$ws.Range("G25").Value = 0.819627
$ws.Range("H25").Value = 2.458881
$ws.Range("I25").Value = 0.05130874782515448
$ws.Range("J25").Value = 0.06184680737567548
$ws.Range("M25").Value = 18.627865
$ws.Range("N25").Value = 37.25573
$ws.Range("O25").Value = 0.07787983970082285
$ws.Range("P25").Value = 0.05330365312071852
$ws.Range("Q25").Value = 15.267901106355
$ws.Range("R25").Value = 91.60740663813
$ws.Range("S25").Value = 0.003995917055872974
$ws.Range("T25").Value = 0.003296660766976902
$ws.Range("G26").Value = 0.819627
$ws.Range("H26").Value = 2.458881
$ws.Range("I26").Value = 0.05130874782515448
$ws.Range("J26").Value = 0.06184680737567548
$ws.Range("M26").Value = 15.028766
$ws.Range("N26").Value = 45.086298
$ws.Range("O26").Value = 0.06283263739463307
$ws.Range("P26").Value = 0.06450724194880479
$ws.Range("Q26").Value = 12.317982390282
$ws.Range("R26").Value = 110.861841512538
$ws.Range("S26").Value = 0.0032238639472706
$ws.Range("T26").Value = 0.003989566967143823
